$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.454701162732887
$ws.Range("B2").Value = -4.796697199710548

$ws.Range("A3").Value = -0.5354128534328406
$ws.Range("B3").Value = 0.7871771828006674

$ws.Range("A4").Value = 0.8480562870409967
$ws.Range("B4").Value = -2.641796420480595

$ws.Range("A5").Value = 0.7314811019627357
$ws.Range("B5").Value = 0.6500669003362431

$ws.Range("A6").Value = -0.8064437839643699
$ws.Range("B6").Value = -2.271224825824301

$ws.Range("A7").Value = -0.0761517268177867
$ws.Range("B7").Value = -0.5544581549475408

$ws.Range("A8").Value = 0.7851717713428624
$ws.Range("B8").Value = 0.7866806101569737

$ws.Range("A9").Value = 0.321178934990406
$ws.Range("B9").Value = -0.09485569883604175
